# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 19-20),
# pushing the existing historical rows down by two (rows 19-116 -> 21-118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("19:20").Insert()

# New row 19
$ws.Cells.Item(19, 1).Value = 10
$ws.Cells.Item(19, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value = "La Araucanía"
$ws.Cells.Item(19, 4).Value = 45243
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(19, 6).Value = 300000000
$ws.Cells.Item(19, 7).Value = "Espárragos"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 480
$ws.Cells.Item(19, 11).Value = 1600
$ws.Cells.Item(19, 12).Value = 1600
$ws.Cells.Item(19, 13).Value = 1600
$ws.Cells.Item(19, 14).Value = "$/kilo"
$ws.Cells.Item(19, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(19, 16).Value = 1600
$ws.Cells.Item(19, 17).Value = 1
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# New row 20
$ws.Cells.Item(20, 1).Value = 10
$ws.Cells.Item(20, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value = "La Araucanía"
$ws.Cells.Item(20, 4).Value = 45243
$ws.Cells.Item(20, 5).Value = 9
$ws.Cells.Item(20, 6).Value = 300000000
$ws.Cells.Item(20, 7).Value = "Espárragos"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 800
$ws.Cells.Item(20, 11).Value = 1600
$ws.Cells.Item(20, 12).Value = 1800
$ws.Cells.Item(20, 13).Value = 1725
$ws.Cells.Item(20, 14).Value = "$/kilo"
$ws.Cells.Item(20, 15).Value = "Región del Maule"
$ws.Cells.Item(20, 16).Value = 1725
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"
